$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " October 13 2020"
$ws.Range("C4").Value = "Super Kings won by 20 runs"
$ws.Range("D4").Value = "Chennai Super Kings"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Dwayne Bravo "
$ws.Range("G4").Value = "'0"
$ws.Range("H4").Value = "'1"
$ws.Range("I4").Value = "'0"
$ws.Range("J4").Value = "'0"
$ws.Range("K4").Value = "'0.00"
$ws.Range("G4:K4").ClearFormats()

$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 10 2020"
$ws.Range("C5").Value = "RCB won by 37 runs"
$ws.Range("D5").Value = "Chennai Super Kings"
$ws.Range("E5").Value = "Royal Challengers Bangalore"
$ws.Range("F5").Value = "Dwayne Bravo "
$ws.Range("G5").Value = "'7"
$ws.Range("H5").Value = "'5"
$ws.Range("I5").Value = "'0"
$ws.Range("J5").Value = "'0"
$ws.Range("K5").Value = "'140.00"
$ws.Range("G5:K5").ClearFormats()
